# The invoice template had a redundant floating "Text Box 2" shape anchored
# in the "TO:" paragraph that duplicated the Date/Ref/Cust SO No/Terms text
# which already exists as plain paragraphs further down in the body. Remove
# that stray floating text box (and its VML fallback) entirely, leaving the
# surrounding runs/paragraphs untouched.

$d = $word.ActiveDocument

for ($i = $d.Shapes.Count; $i -ge 1; $i--) {
    $shape = $d.Shapes.Item($i)
    if ($shape.Name -eq "Text Box 2") {
        $shape.Delete()
    }
}
